# Add new "loose" game worksheet, after "progressive", mirroring the
# layout of the existing slot-machine analysis sheets.

$wb = $excel.ActiveWorkbook

# --- create the new sheet as the last tab, named "loose" -----------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "loose"

# --- column widths (approximate the bestFit widths of the sibling sheets) -
$ws.Columns.Item(1).ColumnWidth = 31.451822916666668
$ws.Columns.Item(2).ColumnWidth = 9.451822916666666
$ws.Columns.Item(3).ColumnWidth = 9.451822916666666
$ws.Columns.Item(4).ColumnWidth = 11.307291666666666

# --- new shared strings must be introduced in this exact order so the ----
# workbook's shared string table lines up with the source (Heart, Gold Bar,
# Horseshoe, Heart+Any+Any, Heart+Heart+Any, Gold Bar+Any+Any,
# Gold Bar+Gold Bar+Any, Gold Bar+Gold Bar+Gold Bar, Horseshoe x3)
$ws.Range("A2").Value = "Heart"
$ws.Range("A6").Value = "Gold Bar"
$ws.Range("A4").Value = "Horseshoe"
$ws.Range("A11").Value = "Heart+Any+Any"
$ws.Range("A12").Value = "Heart+Heart+Any"
$ws.Range("A16").Value = "Gold Bar+Any+Any"
$ws.Range("A17").Value = "Gold Bar+Gold Bar+Any"
$ws.Range("A18").Value = "Gold Bar+Gold Bar+Gold Bar"
$ws.Range("A14").Value = "Horseshoe+Horseshoe+Horseshoe"

# --- remaining labels (reuse strings already present in the workbook) ----
$ws.Range("A3").Value = "Bell"
$ws.Range("A5").Value = "Seven"
$ws.Range("A8").Value = "Combinations"
$ws.Range("A13").Value = "Bell+Bell+Bell"
$ws.Range("A15").Value = "Seven+Seven+Seven"
$ws.Range("A20").Value = "Total"

$ws.Range("B1").Value = "Slot 1"
$ws.Range("C1").Value = "Slot 2"
$ws.Range("D1").Value = "Slot 3"

$ws.Range("B10").Value = "Number"
$ws.Range("C10").Value = "Frequency"
$ws.Range("D10").Value = "Payout"
$ws.Range("E10").Value = "Payout Ratio"

# --- symbol table (rows 2-6): Slot1 / Slot2 / Slot3 reel counts ----------
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = 20

$ws.Range("B3").Value = 6
$ws.Range("C3").Value = 12
$ws.Range("D3").Value = 12

$ws.Range("B4").Value = 10
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 8

$ws.Range("B5").Value = 10
$ws.Range("C5").Value = 6
$ws.Range("D5").Value = 3

$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 2

# --- total combinations ----------------------------------------------------
$ws.Range("B8").Formula = "=SUM(B2:B7)*SUM(C2:C7)*SUM(D2:D7)"

# --- payout table (rows 11-18): combination counts and fixed payouts -----
$ws.Range("B11").Formula = "=B2*SUM(C3:C7)*SUM(D2:D7)"
$ws.Range("D11").Value = 2

$ws.Range("B12").Formula = "=B2*C2*SUM(D2:D7)"
$ws.Range("D12").Value = 4

$ws.Range("B13").Formula = "=B3*C3*D3"
$ws.Range("D13").Value = 5

$ws.Range("B14").Formula = "=B4*C4*D4"
$ws.Range("D14").Value = 6

$ws.Range("B15").Formula = "=B5*C5*D5"
$ws.Range("D15").Value = 10

$ws.Range("B16").Formula = "=B6*SUM(C2:C5)*SUM(D2:D6)"
$ws.Range("D16").Value = 10

$ws.Range("B17").Formula = "=B6*C6*SUM(D2:D5)"
$ws.Range("D17").Value = 25

$ws.Range("B18").Formula = "=B6*C6*D6"
$ws.Range("D18").Value = 1000

# ratio / expected-value columns, written in the same interleaved order as
# the source workbook so the shared-formula group indices (si) line up:
# C12:C16 -> si 0, E12:E16 -> si 1, C17:C18 -> si 2, E17:E18 -> si 3
$ws.Range("C11").Formula = "=B11/`$B`$8"
$ws.Range("C12:C16").Formula = "=B12/`$B`$8"
$ws.Range("E11").Formula = "=C11*D11"
$ws.Range("E12:E16").Formula = "=C12*D12"
$ws.Range("C17:C18").Formula = "=B17/`$B`$8"
$ws.Range("E17:E18").Formula = "=C17*D17"

$ws.Range("C11:C18").NumberFormat = "0.0000"
$ws.Range("E11:E19").NumberFormat = "0.000"

$ws.Range("C20").Formula = "=SUM(C11:C19)"
$ws.Range("C20").NumberFormat = "0.000"
$ws.Range("E20").Formula = "=SUM(E11:E19)"
$ws.Range("E20").NumberFormat = "0.000"

# --- selection matches the source file ------------------------------------
$ws.Range("A7").Select()
